$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update values in column B (Voltaje) for rows 2-5
$ws.Range("B2").Value = 5.6192789999999997
$ws.Range("B3").Value = 5.889278
$ws.Range("B4").Value = 6.5736559999999997
$ws.Range("B5").Value = 6.9976839999999996

# Move the active selection from D5 to D4
$ws.Range("D4").Select()
